# Update input data files to consolidate Cargo Type & Vehicle Type
#
# The "VBDR" sheet (Vehicle Buyer Discount Rate) previously listed one
# discount-rate row per vehicle mode (LDVs, HDVs, aircraft, rail, ships,
# motorbikes). It is being split so each mode has a separate "passenger"
# and "freight" row, and the header row is restructured into a
# "Unit: <label>" cell (A1, italic) + a short column title (B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VBDR")

# --- Relabel the existing "passenger" rows (2-7) ---
$ws.Range("A2").Value = "passenger LDVs"
$ws.Range("A3").Value = "passenger HDVs"
$ws.Range("A4").Value = "passenger aircraft"
$ws.Range("A5").Value = "passenger rail"
$ws.Range("A6").Value = "passenger ships"
$ws.Range("A7").Value = "passenger motorbikes"

# --- Add the new "freight" rows (8-13), same discount rate as the rest ---
$ws.Range("A8").Value = "freight LDVs"
$ws.Range("B8").Value = 0.07

$ws.Range("A9").Value = "freight HDVs"
$ws.Range("B9").Value = 0.07

$ws.Range("A10").Value = "freight aircraft"
$ws.Range("B10").Value = 0.07

$ws.Range("A11").Value = "freight rail"
$ws.Range("B11").Value = 0.07

$ws.Range("A12").Value = "freight ships"
$ws.Range("B12").Value = 0.07

$ws.Range("A13").Value = "freight motorbikes"
$ws.Range("B13").Value = 0.07

# --- Header row: split "Discount Rate (dimensionless)" into a
#     unit label (A1, italic) and a short title (B1) ---
$ws.Range("A1").Value = "Unit: dimensionless"
$ws.Range("A1").Font.Italic = $true
$ws.Range("B1").Value = "Discount Rate"

# --- Widen column A so the longer "passenger"/"freight" labels fit ---
$ws.Columns.Item(1).ColumnWidth = 21.8
